$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"0.4940636666666666"
$ws.Range("H2").Value = [double]"1.482191"
$ws.Range("I2").Value = [double]"0.1416906061387336"
$ws.Range("J2").Value = [double]"0.1416906061387335"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"2.679174"
$ws.Range("N2").Value = [double]"8.037521999999999"
$ws.Range("O2").Value = [double]"0.02942326717729479"
$ws.Range("P2").Value = [double]"0.02942326717729479"
$ws.Range("Q2").Value = [double]"1.323682530078"
$ws.Range("R2").Value = [double]"11.913142770702"
$ws.Range("S2").Value = [double]"0.004169000560932804"
$ws.Range("T2").Value = [double]"0.004169000560932803"

$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"0.4940636666666666"
$ws.Range("H3").Value = [double]"1.482191"
$ws.Range("I3").Value = [double]"0.1416906061387336"
$ws.Range("J3").Value = [double]"0.1416906061387335"
$ws.Range("O3").Value = [double]"0.2465847468531156"
$ws.Range("P3").Value = [double]"0.2465847468531155"
$ws.Range("Q3").Value = [double]"11.09325893777866"
$ws.Range("R3").Value = [double]"99.83933044000797"
$ws.Range("S3").Value = [double]"0.03493874224618412"
$ws.Range("T3").Value = [double]"0.03493874224618411"

$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"0.4940636666666666"
$ws.Range("H4").Value = [double]"1.482191"
$ws.Range("I4").Value = [double]"0.1416906061387336"
$ws.Range("J4").Value = [double]"0.1416906061387335"
$ws.Range("M4").Value = [double]"65.67046766666668"
$ws.Range("N4").Value = [double]"197.011403"
$ws.Range("O4").Value = [double]"0.7212072511207682"
$ws.Range("P4").Value = [double]"0.7212072511207681"
$ws.Range("Q4").Value = [double]"32.44539204710811"
$ws.Range("R4").Value = [double]"292.008528423973"
$ws.Range("S4").Value = [double]"0.1021882925629515"
$ws.Range("T4").Value = [double]"0.1021882925629514"

$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"0.4940636666666666"
$ws.Range("H5").Value = [double]"1.482191"
$ws.Range("I5").Value = [double]"0.1416906061387336"
$ws.Range("J5").Value = [double]"0.1416906061387335"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"0.2535676666666667"
$ws.Range("N5").Value = [double]"0.760703"
$ws.Range("O5").Value = [double]"0.002784734848821526"
$ws.Range("P5").Value = [double]"0.002784734848821525"
$ws.Range("Q5").Value = [double]"0.1252785711414444"
$ws.Range("R5").Value = [double]"1.127507140273"
$ws.Range("S5").Value = [double]"0.0003945707686651766"
$ws.Range("T5").Value = [double]"0.0003945707686651764"

$ws.Range("H6").Value = [double]"6.480663"
$ws.Range("I6").Value = [double]"0.6195214170446748"
$ws.Range("J6").Value = [double]"0.6195214170446747"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"2.679174"
$ws.Range("N6").Value = [double]"8.037521999999999"
$ws.Range("O6").Value = [double]"0.02942326717729479"
$ws.Range("P6").Value = [double]"0.02942326717729479"
$ws.Range("Q6").Value = [double]"5.787607937453999"
$ws.Range("R6").Value = [double]"52.08847143708599"
$ws.Range("S6").Value = [double]"0.01822834417576174"
$ws.Range("T6").Value = [double]"0.01822834417576174"

$ws.Range("H7").Value = [double]"6.480663"
$ws.Range("I7").Value = [double]"0.6195214170446748"
$ws.Range("J7").Value = [double]"0.6195214170446747"
$ws.Range("O7").Value = [double]"0.2465847468531156"
$ws.Range("P7").Value = [double]"0.2465847468531155"
$ws.Range("Q7").Value = [double]"48.503649494216"
$ws.Range("R7").Value = [double]"436.5328454479439"
$ws.Range("S7").Value = [double]"0.1527645317920446"
$ws.Range("T7").Value = [double]"0.1527645317920445"

$ws.Range("H8").Value = [double]"6.480663"
$ws.Range("I8").Value = [double]"0.6195214170446748"
$ws.Range("J8").Value = [double]"0.6195214170446747"
$ws.Range("M8").Value = [double]"65.67046766666668"
$ws.Range("N8").Value = [double]"197.011403"
$ws.Range("O8").Value = [double]"0.7212072511207682"
$ws.Range("P8").Value = [double]"0.7212072511207681"
$ws.Range("Q8").Value = [double]"141.8627233333544"
$ws.Range("R8").Value = [double]"1276.764510000189"
$ws.Range("S8").Value = [double]"0.4468033381972329"
$ws.Range("T8").Value = [double]"0.4468033381972328"

$ws.Range("H9").Value = [double]"6.480663"
$ws.Range("I9").Value = [double]"0.6195214170446748"
$ws.Range("J9").Value = [double]"0.6195214170446747"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"0.2535676666666667"
$ws.Range("N9").Value = [double]"0.760703"
$ws.Range("O9").Value = [double]"0.002784734848821526"
$ws.Range("P9").Value = [double]"0.002784734848821525"
$ws.Range("Q9").Value = [double]"0.5477621984543334"
$ws.Range("R9").Value = [double]"4.929859786089"
$ws.Range("S9").Value = [double]"0.0017252028796356"
$ws.Range("T9").Value = [double]"0.001725202879635599"

$ws.Range("G10").Value = [double]"0.7459539999999999"
$ws.Range("H10").Value = [double]"2.237862"
$ws.Range("I10").Value = [double]"0.2139292596128559"
$ws.Range("J10").Value = [double]"0.2139292596128559"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"2.679174"
$ws.Range("N10").Value = [double]"8.037521999999999"
$ws.Range("O10").Value = [double]"0.02942326717729479"
$ws.Range("P10").Value = [double]"0.02942326717729479"
$ws.Range("Q10").Value = [double]"1.998540561995999"
$ws.Range("R10").Value = [double]"17.986865057964"
$ws.Range("S10").Value = [double]"0.006294497762629921"
$ws.Range("T10").Value = [double]"0.006294497762629921"

$ws.Range("G11").Value = [double]"0.7459539999999999"
$ws.Range("H11").Value = [double]"2.237862"
$ws.Range("I11").Value = [double]"0.2139292596128559"
$ws.Range("J11").Value = [double]"0.2139292596128559"
$ws.Range("O11").Value = [double]"0.2465847468531156"
$ws.Range("P11").Value = [double]"0.2465847468531155"
$ws.Range("Q11").Value = [double]"16.748976773584"
$ws.Range("R11").Value = [double]"150.740790962256"
$ws.Range("S11").Value = [double]"0.05275169232611052"
$ws.Range("T11").Value = [double]"0.05275169232611052"

$ws.Range("G12").Value = [double]"0.7459539999999999"
$ws.Range("H12").Value = [double]"2.237862"
$ws.Range("I12").Value = [double]"0.2139292596128559"
$ws.Range("J12").Value = [double]"0.2139292596128559"
$ws.Range("M12").Value = [double]"65.67046766666668"
$ws.Range("N12").Value = [double]"197.011403"
$ws.Range("O12").Value = [double]"0.7212072511207682"
$ws.Range("P12").Value = [double]"0.7212072511207681"
$ws.Range("Q12").Value = [double]"48.98714803782067"
$ws.Range("R12").Value = [double]"440.884332340386"
$ws.Range("S12").Value = [double]"0.154287333259689"
$ws.Range("T12").Value = [double]"0.154287333259689"

$ws.Range("G13").Value = [double]"0.7459539999999999"
$ws.Range("H13").Value = [double]"2.237862"
$ws.Range("I13").Value = [double]"0.2139292596128559"
$ws.Range("J13").Value = [double]"0.2139292596128559"
$ws.Range("K13").Value = [double]"3"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"0.2535676666666667"
$ws.Range("N13").Value = [double]"0.760703"
$ws.Range("O13").Value = [double]"0.002784734848821526"
$ws.Range("P13").Value = [double]"0.002784734848821525"
$ws.Range("Q13").Value = [double]"0.1891498152206667"
$ws.Range("R13").Value = [double]"1.702348336986"
$ws.Range("S13").Value = [double]"0.0005957362644265073"
$ws.Range("T13").Value = [double]"0.0005957362644265072"

$ws.Range("E14").Value = [double]"1"
$ws.Range("F14").Value = [double]"0.3333333333333333"
$ws.Range("G14").Value = [double]"0.08668033333333335"
$ws.Range("H14").Value = [double]"0.260041"
$ws.Range("I14").Value = [double]"0.02485871720373584"
$ws.Range("J14").Value = [double]"0.02485871720373583"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"2.679174"
$ws.Range("N14").Value = [double]"8.037521999999999"
$ws.Range("O14").Value = [double]"0.02942326717729479"
$ws.Range("P14").Value = [double]"0.02942326717729479"
$ws.Range("Q14").Value = [double]"0.232231695378"
$ws.Range("R14").Value = [double]"2.090085258402"
$ws.Range("S14").Value = [double]"0.000731424677970334"
$ws.Range("T14").Value = [double]"0.0007314246779703339"

$ws.Range("E15").Value = [double]"1"
$ws.Range("F15").Value = [double]"0.3333333333333333"
$ws.Range("G15").Value = [double]"0.08668033333333335"
$ws.Range("H15").Value = [double]"0.260041"
$ws.Range("I15").Value = [double]"0.02485871720373584"
$ws.Range("J15").Value = [double]"0.02485871720373583"
$ws.Range("O15").Value = [double]"0.2465847468531156"
$ws.Range("P15").Value = [double]"0.2465847468531155"
$ws.Range("Q15").Value = [double]"1.946241845645333"
$ws.Range("R15").Value = [double]"17.516176610808"
$ws.Range("S15").Value = [double]"0.00612978048877639"
$ws.Range("T15").Value = [double]"0.006129780488776389"

$ws.Range("E16").Value = [double]"1"
$ws.Range("F16").Value = [double]"0.3333333333333333"
$ws.Range("G16").Value = [double]"0.08668033333333335"
$ws.Range("H16").Value = [double]"0.260041"
$ws.Range("I16").Value = [double]"0.02485871720373584"
$ws.Range("J16").Value = [double]"0.02485871720373583"
$ws.Range("M16").Value = [double]"65.67046766666668"
$ws.Range("N16").Value = [double]"197.011403"
$ws.Range("O16").Value = [double]"0.7212072511207682"
$ws.Range("P16").Value = [double]"0.7212072511207681"
$ws.Range("Q16").Value = [double]"5.692338027502558"
$ws.Range("R16").Value = [double]"51.23104224752301"
$ws.Range("S16").Value = [double]"0.01792828710089487"
$ws.Range("T16").Value = [double]"0.01792828710089487"

$ws.Range("E17").Value = [double]"1"
$ws.Range("F17").Value = [double]"0.3333333333333333"
$ws.Range("G17").Value = [double]"0.08668033333333335"
$ws.Range("H17").Value = [double]"0.260041"
$ws.Range("I17").Value = [double]"0.02485871720373584"
$ws.Range("J17").Value = [double]"0.02485871720373583"
$ws.Range("K17").Value = [double]"3"
$ws.Range("L17").Value = [double]"1"
$ws.Range("M17").Value = [double]"0.2535676666666667"
$ws.Range("N17").Value = [double]"0.760703"
$ws.Range("O17").Value = [double]"0.002784734848821526"
$ws.Range("P17").Value = [double]"0.002784734848821525"
$ws.Range("Q17").Value = [double]"0.02197932986922223"
$ws.Range("R17").Value = [double]"0.197813968823"
$ws.Range("S17").Value = [double]"6.922493609424238E-05"
$ws.Range("T17").Value = [double]"6.922493609424235E-05"
